$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(1.02, 1.085048644128605, 1.087980782939046, 1.087926566658925, 1.099669929262688, 1.068117173912336, 1.089905424368914, 1.09063464894122, 1.09058057227361, 1.102294040925225, 1.091453215085418),
    @(1.02, 1.08643571763501, 1.08913806745221, 1.089182942581487, 1.100996959608187, 1.068650053587837, 1.09095354479349, 1.091610787050656, 1.091655555005366, 1.103441682849266, 1.092502823961218),
    @(1.02, 1.087332004646575, 1.089885779368614, 1.089994488272621, 1.101854623681957, 1.06899293960828, 1.09162998715082, 1.092240698853189, 1.092349161241834, 1.104182714338881, 1.093180226944284),
    @(1.02, 1.087708511926824, 1.090199851368598, 1.09033532858016, 1.102214947923218, 1.069136631032567, 1.091913946078946, 1.092505107142036, 1.092640283998206, 1.104493872731965, 1.093464589126659),
    @(1.02, 1.087771712208968, 1.090252569954156, 1.090392537753033, 1.102275434148457, 1.069160730677033, 1.091961599678299, 1.092549478725027, 1.092689137402578, 1.104546095956155, 1.093512310399589),
    @(1.02, 1.087337036694772, 1.089889977054521, 1.089999043901818, 1.101859439280206, 1.068994861416152, 1.09163378305882, 1.092244233482632, 1.0923530530767, 1.104186873505685, 1.093184028242909),
    @(1.02, 1.085517672735602, 1.088372128798516, 1.088351460425531, 1.10011861811626, 1.068297662695122, 1.090260008928863, 1.090964896804884, 1.090944280349007, 1.102682219601145, 1.091808303196112),
    @(1.02, 1.082301953678646, 1.085688659064118, 1.085437161479364, 1.09704309022976, 1.067054273034987, 1.087825557963856, 1.088697233245837, 1.08844647647644, 1.100018603220093, 1.0893703950315),
    @(1.02, 1.080151218646692, 1.083893485073913, 1.083486558536975, 1.094987067075027, 1.06621521581723, 1.086193117139393, 1.087176262321768, 1.086770652723627, 1.098234354614831, 1.087735635953822),
    @(1.02, 1.079218205109822, 1.08311462737454, 1.082640021225304, 1.094095372269268, 1.065849456870624, 1.085483944104254, 1.086515428573439, 1.08604241832261, 1.09745967398313, 1.087025455811631),
    @(1.02, 1.078871375320308, 1.082825089134557, 1.082325286029999, 1.093763937266375, 1.065713227800211, 1.085220172448212, 1.086269623592112, 1.085771524585224, 1.097171603606538, 1.086761309569582),
    @(1.02, 1.078945783682675, 1.082887206823435, 1.082392811179452, 1.093835041245435, 1.065742466187297, 1.085276768433888, 1.086322365145458, 1.085829650180312, 1.097233410228701, 1.086817985928049),
    @(1.02, 1.079189541538666, 1.083090698893517, 1.082614011139446, 1.0940679802684, 1.065838203697173, 1.085462147893449, 1.086495117260031, 1.086020034253827, 1.097435868546201, 1.087003628647706),
    @(1.02, 1.079339693207604, 1.083216045672927, 1.082750260655017, 1.094211472375255, 1.065897141627185, 1.085576319331625, 1.0866015100875, 1.086137283645432, 1.097560567292446, 1.087117962222421),
    @(1.02, 1.080213102618754, 1.083945142506735, 1.083542699599495, 1.095046215423509, 1.066239438328681, 1.08624013337938, 1.087220072002269, 1.086818928154407, 1.098285723121114, 1.087782718962263),
    @(1.02, 1.080760500919259, 1.084402070984571, 1.084039259067413, 1.09556944240613, 1.066453496139866, 1.086655902729694, 1.087607475523347, 1.087245807880043, 1.098740030918336, 1.088199078752744),
    @(1.02, 1.081079622535084, 1.084668441828436, 1.084328709605573, 1.09587449494345, 1.066578117067019, 1.086898190920435, 1.087833225328268, 1.087494549931139, 1.099004819773556, 1.088441711020507),
    @(1.02, 1.081188406718329, 1.084759242485953, 1.084427373698939, 1.095978486934635, 1.066620569752246, 1.086980767262415, 1.087910163667716, 1.087579322365319, 1.099095071973065, 1.088524404630364),
    @(1.02, 1.080701787541275, 1.084353062209942, 1.083986002033934, 1.095513319308988, 1.06643055413135, 1.086611317723057, 1.087565933182135, 1.087200033592599, 1.098691308843989, 1.088154430430284),
    @(1.02, 1.079117768352619, 1.083030782101625, 1.082548881416844, 1.093999391649514, 1.065810021622966, 1.085407568056703, 1.086444255545531, 1.085963981840723, 1.097376258482819, 1.086948971301331),
    @(1.02, 1.078120284256394, 1.082198045222038, 1.081643605240507, 1.093046250248712, 1.065417726466453, 1.084648676596193, 1.085737031711144, 1.085184538529897, 1.096547582806059, 1.086189002127902),
    @(1.02, 1.078649218544911, 1.082639626080888, 1.08212367258406, 1.093551651415938, 1.065625893598054, 1.085051175076702, 1.086112133818153, 1.085597955009128, 1.096987056829182, 1.08659207220242),
    @(1.02, 1.080728318103635, 1.084375207622021, 1.084010067163685, 1.09553867934344, 1.066440921364522, 1.086631464457202, 1.087584705045396, 1.087220717793817, 1.098713324873464, 1.088174605775102),
    @(1.02, 1.083134488069407, 1.086383472965325, 1.086191918468698, 1.09783916489017, 1.067377493219344, 1.088456571607262, 1.089285081142773, 1.08909406936382, 1.100708689804984, 1.090002304786658)
)

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range($cols[$j] + $rowNum).Value = $rowValues[$j]
    }
}